# "fix xcel header names"
# Append "(%)" to the growth/surprise header labels on both sheets, and
# widen the columns that hold those (now-longer) headers to match.

$wb = $excel.ActiveWorkbook

# --- IncomeReport sheet -----------------------------------------------
$wsIncome = $wb.Worksheets.Item("IncomeReport")

$wsIncome.Range("C1").Value = "Rev. Growth (%)"
$wsIncome.Range("E1").Value = "EPS Growth (%)"
$wsIncome.Range("G1").Value = "FCF Growth (%)"

# Widen columns C, E, G by 4 characters to fit the longer headers.
$wsIncome.Columns.Item(3).ColumnWidth = 14.877604166666666
$wsIncome.Columns.Item(5).ColumnWidth = 13.877604166666666
$wsIncome.Columns.Item(7).ColumnWidth = 13.877604166666666

# --- EarningsReport sheet ----------------------------------------------
$wsEarnings = $wb.Worksheets.Item("EarningsReport")

$wsEarnings.Range("D1").Value = "EPS Surprise (%)"
$wsEarnings.Range("E1").Value = "EPS Growth Quarter and Year Forecast (%)"
$wsEarnings.Range("F1").Value = "Revenue Growth Quarter and Year Forecast (%)"

# Widen columns D, E, F by 4 characters to fit the longer headers.
$wsEarnings.Columns.Item(4).ColumnWidth = 15.877604166666666
$wsEarnings.Columns.Item(5).ColumnWidth = 39.877604166666664
$wsEarnings.Columns.Item(6).ColumnWidth = 43.877604166666664
